$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial": update row 7 (Hernández Mendoza Delfina / 6APV) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E7").Value = 14
$ws2.Range("F7").Value = 15
$ws2.Range("G7").Value = 48.28
$ws2.Range("H7").Value = 51.72
$ws2.Range("I7").Value = 8.800000000000001
$ws2.Range("J7").Value = 15
$ws2.Range("K7").Value = 51.72

# --- Sheet "3er Parcial": update I7 only ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("I7").Value = 8.800000000000001
